# GitHub Actions "Updated cryptos list" refresh: rewrites the Price (col D)
# and Volume(1h) (col E) text for every coin row, and for a handful of rows
# the ranking shuffled so Coin/Link/Price/Volume (cols B-E) moved to a
# neighboring row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All of column D/E (and the B/C coin name/link on a few rows) are stored as
# plain text in this sheet, even values that look numeric (e.g. "214.05",
# "1.004", or multi-dot "25.837.68" style prices). Assigning those bare to
# .Value would let Excel auto-convert them into real numbers, so any cell
# whose new text parses as a plain number is written with a leading
# apostrophe to force Text, then its style is reset to Normal so the sheet
# doesn't end up with a stray quote-prefix style applied.

$ws.Range("D2").Value = '25.837.68'
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").Value = '1.635.24'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("E4").Value = '  -0.69%  '
$ws.Range("D5").Value = '''214.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("D6").Value = '''0.5017'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("D7").Value = '''1.002'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.74%  '
$ws.Range("D8").Value = '''0.2555'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.78%  '
$ws.Range("E9").Value = '  -0.74%  '
$ws.Range("D10").Value = '''19.34'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.92%  '
$ws.Range("D11").Value = '''0.07780'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("D12").Value = '1.647.68'
$ws.Range("E12").Value = '  +0.66%  '
$ws.Range("D13").Value = '''4.235'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.42%  '
$ws.Range("D14").Value = '1.863.11'
$ws.Range("E14").Value = '  +0.09%  '
$ws.Range("E15").Value = '  -0.88%  '
$ws.Range("D16").Value = '0.0₅7857'
$ws.Range("E16").Value = '  -1.08%  '
$ws.Range("D17").Value = '''64.14'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.14%  '
$ws.Range("D18").Value = '25.878.49'
$ws.Range("E18").Value = '  -0.35%  '
$ws.Range("D19").Value = '''1.003'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.70%  '
$ws.Range("D20").Value = '''194.91'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.79%  '
$ws.Range("E21").Value = '  +0.97%  '
$ws.Range("D22").Value = '''9.853'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.35%  '
$ws.Range("D23").Value = '''5.941'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.35%  '
$ws.Range("D24").Value = '''1.004'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.66%  '
$ws.Range("D25").Value = '''1.892'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.49%  '
$ws.Range("D26").Value = '''139.65'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.19%  '
$ws.Range("D27").Value = '''0.1127'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.98%  '
$ws.Range("D28").Value = '''6.780'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.28%  '
$ws.Range("D29").Value = '''15.60'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.91%  '
$ws.Range("D30").Value = '''1.235'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("D31").Value = '''0.04841'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.14%  '
$ws.Range("D32").Value = '''3.231'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.96%  '
$ws.Range("E33").Value = '  -1.13%  '
$ws.Range("D34").Value = '''1.522'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.11%  '
$ws.Range("D35").Value = '''2.360'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.12%  '
$ws.Range("D36").Value = '''2.595'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.55%  '
$ws.Range("D37").Value = '''0.8816'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.69%  '
$ws.Range("D38").Value = '1.124.02'
$ws.Range("E38").Value = '  +0.53%  '
$ws.Range("D39").Value = '''0.5494'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.88%  '
$ws.Range("D40").Value = '''0.01556'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.34%  '
$ws.Range("D41").Value = '''1.003'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.65%  '
$ws.Range("D42").Value = '''5.635'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.34%  '
$ws.Range("D43").Value = '''0.8094'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.79%  '
$ws.Range("E44").Value = '  -0.35%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '0.0₈121'
$ws.Range("E45").Value = '  +8.91%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.775.74'
$ws.Range("E46").Value = '  +0.21%  '
$ws.Range("D47").Value = '''0.4515'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.43%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '''54.98'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.52%  '
$ws.Range("B49").Value = 'Frax'
$ws.Range("C49").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D49").Value = '''1.001'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.45%  '
$ws.Range("D50").Value = '''0.05032'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("D51").Value = '''1.006'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.31%  '
